# Flujos de refinanciacion Especial y registro de informe vista verificacion

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reference/code value in A2 changes from 1332236 to 3534375
$ws.Range("A2").Value = "3534375"

# The sheet view scrolled one column to the right (topLeftCell K1 -> L1)
# and the active selection moved from A2 to R2.
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.TopLeftCell = $ws.Range("L1")
$ws.Range("R2").Select()
